$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodePostResults")

$ws.Range("A2").Value = "https://scrapeops.io/python-web-scraping-playbook/python-requests-post-requests/#post-json-data-using-python-requests"
$ws.Range("B2").Value = 405
